$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "Save" in H1, copying the header formatting (bold, border, centered)
# from the existing "sum" header in G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values for each data row (H2:H20)
$saveValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 1
    20 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
